# Create two new slides ("Title and Content" layout) after the existing
# slide and populate them with the Go-testing content.

$p = $ppt.ActivePresentation

# --- Slide 2: "Create & Run Test File" -------------------------------
$s2 = $p.Slides.Add(2, 2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Create & Run Test File"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "To create a test, create a file named *_test.go" + [char]13 + "To run all the tests in a package type " + [char]8220 + "go test" + [char]8221

# Split "test.go" into its own run (mirrors the proofing-tool run break
# PowerPoint inserts around the flagged word).
$body2.Characters(41, 7).Text = "test.go"

# --- Slide 3: "Test Cases" --------------------------------------------
$s3 = $p.Slides.Add(3, 2)

$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Test Cases "

$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "Make sure created deck has the correct number of cards." + [char]13 + "Write if statement to see if the deck has the right number of cards" + [char]13 + "If not then test go test handler that something is wrong."

$body3.Paragraphs(2).IndentLevel = 2
$body3.Paragraphs(3).IndentLevel = 2

# Split the closing sentence into its own run, matching the source deck.
$body3.Paragraphs(3).Characters(49, 9).Text = "is wrong."
